$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.86'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.95'

$ws.Range("B4").Value = 'LEO'
$ws.Range("C4").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '3.500'
$ws.Range("E4").Value = '3LEOLEO'

$ws.Range("B5").Value = 'HuobiToken'
$ws.Range("C5").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '5.011'
$ws.Range("E5").Value = '4HuobiTokenHT'

$ws.Range("B6").Value = 'Cronos'
$ws.Range("C6").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.05612'
$ws.Range("E6").Value = '5CronosCRO'

$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.521'
$ws.Range("E7").Value = '6KuCoinTokenKCS'

$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.988'
$ws.Range("E8").Value = '7GateTokenGT'

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8103'
$ws.Range("E9").Value = '8MXTokenMX'

$ws.Range("B10").Value = 'FTXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8380'
$ws.Range("E10").Value = '9FTXTokenFTT'

$ws.Range("B11").Value = 'One'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.009633'
$ws.Range("E11").Value = '10OneONEBestin24h'

$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1338'
$ws.Range("E12").Value = '11WazirXWRX'

$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03301'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B14").Value = 'MandalaExchangeToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.06954'
$ws.Range("E14").Value = '13MandalaExchangeTokenMDX'

$ws.Range("B15").Value = 'BitrueCoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.02843'
$ws.Range("E15").Value = '14BitrueCoinBTR'

$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.09411'
$ws.Range("E16").Value = '15BitMartTokenBMX'

$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001515'
$ws.Range("E17").Value = '16BitForexTokenBF'

$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006236'
$ws.Range("E18").Value = '17TigerCashTCH'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.092'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3184'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.746'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04682'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004522'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009701'
$ws.Range("E27").Value = '26NitroExNTX'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001940'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03631'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006241'
$ws.Range("E41").Value = '40KickTokenKICK'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1052'
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002723'
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008358'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005267'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.2000'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002107'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
